$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-39 from 45180 to 45181
$ws.Range("C2:C39").Value = 45181
